$wb = $excel.ActiveWorkbook

# --- Monday sheet (sheet2) updates ---
$monday = $wb.Worksheets.Item("Monday")
$monday.Range("D5").Value = "anna"
$monday.Range("E5").Value = "jack"
$monday.Range("F5").Value = "anna"
$monday.Range("F6").Value = "anna"
$monday.Range("A15").Value = 21
$monday.Range("A16").Value = 23
$monday.Range("B16").Value = "anna"
$monday.Range("C16").Value = "jack"
$monday.Range("E16").Value = "anna"
$monday.Range("D6").Select() | Out-Null

# --- pictures sheet (sheet7) updates ---
$pictures = $wb.Worksheets.Item("pictures")
$pictures.Range("B3").Value = 4000
$pictures.Range("B3").Select() | Out-Null

# --- names sheet (sheet1) selection update ---
$names = $wb.Worksheets.Item("names")
$names.Range("A4").Select() | Out-Null

# Make Monday the active sheet/tab last so it ends up active
$monday.Select() | Out-Null
$monday.Range("D6").Select() | Out-Null
